$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1934306569343066
$ws.Range("C2").Value = 0.583941605839416
$ws.Range("J2").Value = 0.01094890510948905
$ws.Range("P2").Value = 0.145985401459854
$ws.Range("S2").Value = 0.06569343065693431
$ws.Range("C3").Value = 0.03067484662576687
$ws.Range("J3").Value = 0.01840490797546012
$ws.Range("P3").Value = 0.7914110429447853
$ws.Range("S3").Value = 0.1595092024539877
$ws.Range("P4").Value = 0.8333333333333334
$ws.Range("S4").Value = 0.1666666666666667
$ws.Range("P5").Value = 1
$ws.Range("B6").Value = 0.06635071090047394
$ws.Range("D6").Value = 0.009478672985781991
$ws.Range("F6").Value = 0.03791469194312796
$ws.Range("J6").Value = 0.2322274881516588
$ws.Range("O6").Value = 0.004739336492890996
$ws.Range("Q6").Value = 0.1895734597156398
$ws.Range("R6").Value = 0.0995260663507109
$ws.Range("S6").Value = 0.3601895734597156
$ws.Range("B7").Value = 0.0963855421686747
$ws.Range("D7").Value = 0.03012048192771084
$ws.Range("F7").Value = 0.05421686746987952
$ws.Range("J7").Value = 0.04819277108433735
$ws.Range("O7").Value = 0.01204819277108434
$ws.Range("Q7").Value = 0.2168674698795181
$ws.Range("R7").Value = 0.1204819277108434
$ws.Range("S7").Value = 0.4216867469879518
$ws.Range("B8").Value = 0.09457092819614711
$ws.Range("D8").Value = 0.01751313485113835
$ws.Range("E8").Value = 0.001751313485113835
$ws.Range("F8").Value = 0.05078809106830123
$ws.Range("J8").Value = 0.09632224168126094
$ws.Range("O8").Value = 0.01576182136602452
$ws.Range("Q8").Value = 0.1751313485113835
$ws.Range("R8").Value = 0.1436077057793345
$ws.Range("S8").Value = 0.404553415061296
$ws.Range("B9").Value = 0.1068376068376068
$ws.Range("D9").Value = 0.02136752136752137
$ws.Range("F9").Value = 0.04700854700854701
$ws.Range("J9").Value = 0.08547008547008547
$ws.Range("O9").Value = 0.004273504273504274
$ws.Range("Q9").Value = 0.2222222222222222
$ws.Range("R9").Value = 0.0811965811965812
$ws.Range("S9").Value = 0.4316239316239316
$ws.Range("B10").Value = 0.08314087759815242
$ws.Range("D10").Value = 0.02463433410315627
$ws.Range("F10").Value = 0.07467282525019246
$ws.Range("J10").Value = 0.09314857582755966
$ws.Range("O10").Value = 0.01077752117013087
$ws.Range("Q10").Value = 0.2517321016166282
$ws.Range("R10").Value = 0.09622786759045419
$ws.Range("S10").Value = 0.365665896843726
$ws.Range("G11").Value = 0.1646090534979424
$ws.Range("J11").Value = 0.08230452674897119
$ws.Range("K11").Value = 0.2345679012345679
$ws.Range("L11").Value = 0.5144032921810699
$ws.Range("S11").Value = 0.00411522633744856
$ws.Range("G12").Value = 0.751937984496124
$ws.Range("J12").Value = 0.1937984496124031
$ws.Range("L12").Value = 0.0310077519379845
$ws.Range("S12").Value = 0.02325581395348837
$ws.Range("G13").Value = 0.7083333333333334
$ws.Range("J13").Value = 0.2708333333333333
$ws.Range("S13").Value = 0.02083333333333333
$ws.Range("F15").Value = 0.02083333333333333
$ws.Range("H15").Value = 0.1875
$ws.Range("I15").Value = 0.06770833333333333
$ws.Range("J15").Value = 0.3541666666666667
$ws.Range("K15").Value = 0.05208333333333334
$ws.Range("M15").Value = 0.02083333333333333
$ws.Range("O15").Value = 0.03645833333333334
$ws.Range("S15").Value = 0.2604166666666667
$ws.Range("F16").Value = 0.009478672985781991
$ws.Range("H16").Value = 0.2274881516587678
$ws.Range("I16").Value = 0.07582938388625593
$ws.Range("J16").Value = 0.4170616113744076
$ws.Range("K16").Value = 0.07582938388625593
$ws.Range("M16").Value = 0.02369668246445497
$ws.Range("O16").Value = 0.04265402843601896
$ws.Range("S16").Value = 0.1279620853080569
$ws.Range("F17").Value = 0.007233273056057866
$ws.Range("H17").Value = 0.1934900542495479
$ws.Range("I17").Value = 0.1157323688969259
$ws.Range("J17").Value = 0.4358047016274864
$ws.Range("K17").Value = 0.0650994575045208
$ws.Range("M17").Value = 0.0216998191681736
$ws.Range("N17").Value = 0.001808318264014466
$ws.Range("O17").Value = 0.02893309222423146
$ws.Range("S17").Value = 0.1301989150090416
$ws.Range("F18").Value = 0.01893939393939394
$ws.Range("H18").Value = 0.2651515151515151
$ws.Range("I18").Value = 0.07954545454545454
$ws.Range("J18").Value = 0.3484848484848485
$ws.Range("K18").Value = 0.07575757575757576
$ws.Range("M18").Value = 0.02651515151515152
$ws.Range("N18").Value = 0.003787878787878788
$ws.Range("O18").Value = 0.05681818181818182
$ws.Range("S18").Value = 0.125
$ws.Range("F19").Value = 0.01515151515151515
$ws.Range("H19").Value = 0.2393939393939394
$ws.Range("I19").Value = 0.0946969696969697
$ws.Range("J19").Value = 0.3840909090909091
$ws.Range("K19").Value = 0.07727272727272727
$ws.Range("M19").Value = 0.01742424242424243
$ws.Range("N19").Value = 0.0007575757575757576
$ws.Range("O19").Value = 0.06363636363636363
$ws.Range("S19").Value = 0.1075757575757576
